$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (B2:O2) with the values shown in the diff
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 2
